# Code Review Comments added
# Adds the new Settings rows (exception-email addresses + SE_/BE_ message
# catalogue + runtime local folder path) that were appended below the
# existing SystemException/BusinessException email templates, and moves
# the sheet's active selection down to the new last populated cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

$rows = @(
    @("ExceptionEmail",          "Vinusangari.Saravanan@omes.ok.gov"),
    @("BusinessExceptionEmail",  "Vinusangari.Saravanan@omes.ok.gov"),
    @("SE_1",                    "Unable to login to Workday Application"),
    @("BE_1",                    "Unable to find JR in Workday " + [char]0x2013 + " ""<JR Number>"""),
    @("BE_2",                    "Invalid Job Transfer Type"),
    @("SE_2",                    "Unable to find Hire Task for the employee in inbox"),
    @("SE_3",                    "Unable to update employee details in hire task in inbox"),
    @("SE_4",                    "Unable to update employee salary details in hire task in inbox"),
    @("SE_5",                    "Unable to submit disposition in workday"),
    @("BE_3",                    "Unable to find DOH value for the employee in sharepoint list"),
    @("SE_6",                    "Unable to find Change Job Task for the employee in inbox"),
    @("RuntimeLocalFolderPath",  "C:\Users\<Username>\Documents\NHC DispositionWorkday")
)

$startRow = 27
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

$ws.Activate()
$ws.Range("B46").Select()

Write-Host "Added" $rows.Count "rows to Settings starting at row" $startRow
